$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.324.22'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.551.92'
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '517.62'
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.13'
$ws.Range("E6").Value = '  -1.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E8").Value = '  -1.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.563.07'
$ws.Range("E9").Value = '  -3.33%  '
$ws.Range("E10").Value = '  -4.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0989'
$ws.Range("E11").Value = '  -3.04%  '
$ws.Range("E12").Value = '  -3.19%  '
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.004.25'
$ws.Range("E14").Value = '  -3.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.308.59'
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.95'
$ws.Range("E16").Value = '  -4.47%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  -2.94%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.562.97'
$ws.Range("E18").Value = '  -3.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '330.21'
$ws.Range("E19").Value = '  -2.29%  '
$ws.Range("E20").Value = '  -2.84%  '
$ws.Range("E21").Value = '  -2.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.17'
$ws.Range("E22").Value = '  -2.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.38'
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.397'
$ws.Range("E27").Value = '  -4.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.663.70'
$ws.Range("E28").Value = '  -4.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.88'
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0735'
$ws.Range("E31").Value = '  -7.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.23'
$ws.Range("E32").Value = '  -6.08%  '
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.49'
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.42'
$ws.Range("E35").Value = '  -2.04%  '
$ws.Range("E36").Value = '  -3.60%  '
$ws.Range("E37").Value = '  -4.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.831'
$ws.Range("E38").Value = '  -7.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.72'
$ws.Range("E39").Value = '  -2.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.820'
$ws.Range("E40").Value = '  -4.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.42'
$ws.Range("E41").Value = '  -2.03%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.46'
$ws.Range("E43").Value = '  -2.68%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '267.19'
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.64'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0947'
$ws.Range("E46").Value = '  -1.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.581'
$ws.Range("E47").Value = '  -4.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.55'
$ws.Range("E48").Value = '  -5.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0515'
$ws.Range("E49").Value = '  -2.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.958.05'
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("E51").Value = '  -3.64%  '
